$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1)
$ws.Range("B1").Value = "sr"
$ws.Range("C1").Value = "nsecode"
$ws.Range("D1").Value = "name"
$ws.Range("E1").Value = "bsecode"
$ws.Range("F1").Value = "per_chg"
$ws.Range("G1").Value = "close"
$ws.Range("H1").Value = "volume"

# Reuse the existing header style from A1 for the new header cells
$ws.Range("A1").Copy()
$ws.Range("B1:H1").PasteSpecial(-4122)

# Data row 2
$ws.Range("A2").Value = "10/06/2024 04:45:40"
$ws.Range("B2").Value = 1
$ws.Range("C2").Value = "MAZDOCK"
$ws.Range("D2").Value = "Mazagon Dock Shipbuilders Ltd"
$ws.Range("E2").Value = "'543237"
$ws.Range("E2").ClearFormats()
$ws.Range("F2").Value = -1.11
$ws.Range("G2").Value = 3118
$ws.Range("H2").Value = 614819

# Data row 3
$ws.Range("A3").Value = "10/06/2024 04:45:40"
$ws.Range("B3").Value = 2
$ws.Range("C3").Value = "NMDC"
$ws.Range("D3").Value = "Nmdc Limited"
$ws.Range("E3").Value = "'526371"
$ws.Range("E3").ClearFormats()
$ws.Range("F3").Value = -0.99
$ws.Range("G3").Value = 255.95
$ws.Range("H3").Value = 2119786
